$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Numeric columns that resolve to 0
$ws.Range("B2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 0

# Text columns that resolve to an empty string; touch a formatting
# property (no-op value) so the empty cell is still persisted in the
# sheet instead of being dropped as a fully blank cell.
$textCols = @("A2","C2","D2","E2","F2","G2","H2","J2","Q2","R2","S2","T2")
foreach ($addr in $textCols) {
    $cell = $ws.Range($addr)
    $cell.Value = ""
    $cell.Font.Bold = $false
}
